$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.788.44"
$ws.Range("E2").Value = "  -0.44%  "

$ws.Range("D3").Value = "3.408.30"
$ws.Range("E3").Value = "  -0.20%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "412.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.87%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.18"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.66%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.620"
$ws.Range("D7").Style = "Normal"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("E9").Value = "  -1.45%  "

$ws.Range("E10").Value = "  -5.94%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.63"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.17%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000219"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.30%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.13"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.39%  "

$ws.Range("D14").Value = "3.947.48"
$ws.Range("E14").Value = "  -0.12%  "

$ws.Range("E15").Value = "  -0.05%  "

$ws.Range("E16").Value = "  -1.86%  "

$ws.Range("D17").Value = "3.398.12"
$ws.Range("E17").Value = "  -0.82%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.45"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.58%  "

$ws.Range("E19").Value = "  +0.75%  "

$ws.Range("D20").Value = "61.853.19"
$ws.Range("E20").Value = "  -0.10%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "479.77"
$ws.Range("D21").Style = "Normal"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "90.62"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.31%  "

$ws.Range("E23").Value = "  +2.59%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.01%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.30"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.24%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.72"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +9.92%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "33.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.08%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.74"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.03%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.80"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.70%  "

$ws.Range("E30").Value = "  +0.11%  "

$ws.Range("E31").Value = "  -3.59%  "

$ws.Range("E32").Value = "  -1.75%  "

$ws.Range("E33").Value = "  -3.52%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "40.83"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.54%  "

$ws.Range("E35").Value = "  -0.70%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "56.81"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.63%  "

$ws.Range("E37").Value = "  -2.85%  "

$ws.Range("E38").Value = "  +0.13%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.98%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "149.14"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.29%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.323"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.67%  "

$ws.Range("E42").Value = "  -0.53%  "

$ws.Range("E43").Value = "  -0.66%  "

$ws.Range("E44").Value = "  +4.25%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.58"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.48%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.19"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.61%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.34"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +18.96%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "16.44"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.23%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.01"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.53%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "112.98"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +13.69%  "

$ws.Range("D51").Value = "0.0₃0510"
$ws.Range("E51").Value = "  +12.24%  "
